$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: "Hours Logged" -> "Minutes Logged" ---
$ws.Range("C1").Value = "Minutes Logged"

# --- Row 2: finish time edited in place, minutes value replaces hours ---
$ws.Range("B2").Value = "2:30PM 6-15-2017"
$ws.Range("C2").Value = 90

# --- Row 3: finish time edited in place, minutes value replaces hours ---
$ws.Range("B3").Value = "9:30AM 10-27-2017"
$ws.Range("C3").Value = 139

# --- Row 4: brand new data row ---
$ws.Range("A4").Value = "10:40PM 10-26-2017"
$ws.Range("B4").Value = "12:12AM 10-27-2017"
$ws.Range("C4").Value = 92

# --- Row 5: new data row, minutes column right aligned ---
$ws.Range("A5").Value = "3:40PM 10-28-2017"
$ws.Range("B5").Value = "6:39PM 10-28-2017"
$ws.Range("C5").Value = 179
$ws.Range("C5").HorizontalAlignment = -4152

# --- Row 6 ---
$ws.Range("A6").Value = "8:00PM 10-28-2017"
$ws.Range("B6").Value = "11:13PM 10-28-2017"
$ws.Range("C6").Value = 193

# --- Row 7 ---
$ws.Range("A7").Value = "8:15PM 11-22-2017"
$ws.Range("B7").Value = "1:15AM 11-23-2017 "
$ws.Range("C7").Value = 300

# --- Row 8: finish time formatted as a time-of-day string ---
$ws.Range("A8").Value = "12:28PM 11-23-2017"
$ws.Range("B8").Value = "01:19PM 11-23-2017"
$ws.Range("B8").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("C8").Value = 51

# --- Row 9 ---
$ws.Range("A9").Value = "4:10PM 11-23-2017"
$ws.Range("B9").Value = "8:06PM 11-23-2017"
$ws.Range("C9").Value = 236

# --- Row 10 ---
$ws.Range("A10").Value = "11:18AM 11-24-2017"
$ws.Range("B10").Value = "1:59PM 11-24-2017"
$ws.Range("C10").Value = 161

# --- Row 17: totals row ---
$ws.Range("A17").Value = "Total Project Hours:"
$ws.Range("C17").Formula = "=SUM(C2:C16)/60"

# --- cosmetic: move the active selection to match the edited range ---
$ws.Range("C10").Select() | Out-Null
